# Generate Report for Handoff
# - Update "Latest Handoff Datetime" / "Latest HO Xliff Generate Date" values
#   for the batch of files that were just handed off (zh-cn batch moved
#   2016-09-06 02:22:44 -> 02:22:59, de-de / Overview batch moved
#   2016-09-06 02:22:49 -> 02:23:10).
# - Mark the Priority column ("ht") for the rows that are "Ready for handoff"
#   and were part of this handoff batch, on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$rows = @(8, 9, 11, 12, 13, 14)

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZh.Range("H$r").Value = "2016-09-06 02:22:59"
    $wsZh.Range("E$r").Value = "ht"
}

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDe.Range("H$r").Value = "2016-09-06 02:23:10"
    $wsDe.Range("E$r").Value = "ht"
}

# --- Overview sheet (Latest HO Xliff Generate Date mirrors the de-de batch time) ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-09-06 02:23:10"
}
